$wb = $excel.ActiveWorkbook

# The "Jill" sheet (sheet2) had a duplicate year-2025 row; delete the
# duplicate row (row 8) which shifts everything below it up by one row.
$wsJill = $wb.Worksheets.Item("Jill")
$wsJill.Rows.Item(8).Delete()

# Make "Jill" the active sheet/tab, with row 8 (now year 2026) selected
# as a whole row, matching the new selection left behind by the edit.
[void]$wsJill.Activate()
[void]$wsJill.Range("A8:XFD8").Select()
